$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# Row 14
$ws.Range("F15").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 6
$ws.Range("K14").Value = 50
$ws.Range("L14").Value = 20
$ws.Range("M14").Value = -57.142857142857
$ws.Range("N14").Value = -70

# Row 15
$ws.Range("C22").Copy($ws.Range("C15"))
$ws.Range("C15").Value = "0"
$ws.Range("C22").Copy($ws.Range("D15"))
$ws.Range("D15").Value = "0"
$ws.Range("E22").Copy($ws.Range("E15"))
$ws.Range("E15").Value = "***.*"
$ws.Range("N15").Value = -64.705882352941

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 22.222222222222
$ws.Range("I16").Value = 156
$ws.Range("J16").Value = 118
$ws.Range("K16").Value = 32.203389830508
$ws.Range("L16").Value = 6.122448979591
$ws.Range("M16").Value = -42.007434944237
$ws.Range("N16").Value = -82.332955832389

# Row 17
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 36
$ws.Range("H17").Value = -13.888888888888
$ws.Range("I17").Value = 362
$ws.Range("J17").Value = 356
$ws.Range("K17").Value = 1.685393258426
$ws.Range("L17").Value = -3.208556149732
$ws.Range("M17").Value = 49.586776859504
$ws.Range("N17").Value = 1.685393258426

# Row 18
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 120
$ws.Range("I18").Value = 184
$ws.Range("J18").Value = 123
$ws.Range("K18").Value = 49.593495934959
$ws.Range("L18").Value = 7.602339181286
$ws.Range("M18").Value = -44.072948328267
$ws.Range("N18").Value = -86.588921282798

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 62.5
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 2.272727272727
$ws.Range("I19").Value = 533
$ws.Range("J19").Value = 387
$ws.Range("K19").Value = 37.726098191214
$ws.Range("L19").Value = 3.495145631067
$ws.Range("M19").Value = 39.528795811518
$ws.Range("N19").Value = 4.921259842519

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 20
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 106.666666666667
$ws.Range("I20").Value = 233
$ws.Range("J20").Value = 124
$ws.Range("K20").Value = 87.903225806451
$ws.Range("L20").Value = -10.03861003861
$ws.Range("M20").Value = -30.03003003003
$ws.Range("N20").Value = -92.009602194787

# Row 21
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 38.709677419354
$ws.Range("F21").Value = 145
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = 23.931623931623
$ws.Range("I21").Value = 1486
$ws.Range("J21").Value = 1149
$ws.Range("K21").Value = 29.329852045256
$ws.Range("L21").Value = -1.131071190951
$ws.Range("M21").Value = -6.246056782334
$ws.Range("N21").Value = -75.595335851535

# Row 24
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 17.857142857142
$ws.Range("F24").Value = 129
$ws.Range("G24").Value = 111
$ws.Range("H24").Value = 16.216216216216
$ws.Range("I24").Value = 1307
$ws.Range("J24").Value = 854
$ws.Range("K24").Value = 53.044496487119
$ws.Range("L24").Value = 51.976744186046
$ws.Range("M24").Value = 82.033426183844

# Row 25
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -30
$ws.Range("F25").Value = 51
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = 13.333333333333
$ws.Range("I25").Value = 569
$ws.Range("J25").Value = 432
$ws.Range("K25").Value = 31.712962962963
$ws.Range("L25").Value = 25.884955752212
$ws.Range("M25").Value = 9.003831417624

# Row 26
$ws.Range("C22").Copy($ws.Range("C26"))
$ws.Range("C26").Value = "0"
$ws.Range("C22").Copy($ws.Range("D26"))
$ws.Range("D26").Value = "0"
$ws.Range("E22").Copy($ws.Range("E26"))
$ws.Range("E26").Value = "***.*"

# Row 27
$ws.Range("C22").Copy($ws.Range("C27"))
$ws.Range("C27").Value = "0"
$ws.Range("F15").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 2
$ws.Range("H15").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 53
$ws.Range("J27").Value = 44
$ws.Range("K27").Value = 20.454545454545
$ws.Range("L27").Value = 17.777777777777

# Row 28
$ws.Range("F15").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("H15").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("J28").Value = 33
$ws.Range("K28").Value = -33.333333333333
$ws.Range("M28").Value = -50

# Row 29
$ws.Range("F15").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("H15").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -33.333333333333
$ws.Range("J29").Value = 27
$ws.Range("K29").Value = -33.333333333333
$ws.Range("M29").Value = -37.931034482758
